$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New note cells (red text), added as part of switching the .csv read-in
# over to map(2): one note attached to the rater/norm_group columns and
# one attached to the CV columns.
# Write F9 first so its string lands in sharedStrings before B8's string,
# matching the author's edit order.
$ws.Range("F9").Value = "Calculate actual intervals, based on these CVs"
$ws.Range("F9").Font.Color = 255

$ws.Range("B8").Value = "Collapse rater and norm_group into single col"
$ws.Range("B8").Font.Color = 255

# Reflect the zoomed-in view and new active selection left by the edit.
$excel.ActiveWindow.Zoom = 159
$ws.Range("B8").Select() | Out-Null
